# "Fin des grilles aleatoire"
# Adds the new journal entry for 25.03.2020 (row 24) describing the
# completion of the random grid-choice construction, and updates the
# active selection to K24 to match the author's saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 24 previously was a blank template row (same style as rows 25-27).
# Pick up the date-cell formatting (centered, bordered, numFmtId=14 date
# format) from the row above so A24 renders as a date like the rest of
# the "Date" column, then overwrite the value.
$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)

$ws.Range("A24").Value = 43915
$ws.Range("B24").Value = "2.25 h "
$ws.Range("C24").Value = "travail pratique"
$ws.Range("D24").Value = "MA-20"
$ws.Range("E24").Value = "fin de la construction des choix aléatoire de grille."
$ws.Range("K24").Value = 7

# Move the active selection to K24 (where the author left off).
$ws.Range("K24").Select()
